$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Title (D) and Description (E) columns with new product copy,
# and fix the tags value in O6.
$ws.Range("D2").Value = "Pretty handmade T-Shirt"
$ws.Range("D3").Value = "Handmade T-Shirt - white"
$ws.Range("D4").Value = "Handmade T-Shirt - blue"
$ws.Range("D5").Value = "Handmade T-Shirt - black"
$ws.Range("D6").Value = "Handmade T-Shirt - green"

$ws.Range("E2:E6").Value = "Handmade high grade item. And so on, and so on…"

$ws.Range("O6").Value = "tag3,tag4,tag5"

# Restore the view to the top-left corner and move the selection to A3,
# matching how the workbook was left after the edits were made.
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A3").Select()
